$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (shared strings changed, and Candidate ID incremented)
$ws.Range("A2").Value = "OOkLX161"
$ws.Range("B2").Value = 23071823
$ws.Range("C2").Value = "urfnboj78"
$ws.Range("D2").Value = "KN&8w7j%"
$ws.Range("F2").Value = "RXpSxSmD"
$ws.Range("G2").Value = "reGn"
